$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the timestamp in A2
$ws.Range("A2").Value = "2020-03-12 20:30:54.168122"

# Labels to insert in column A (row -> label text); the value that used to
# live in column A for that row moves over to column B, unchanged.
$labels = @{
    4  = "Fwd mean draft, m: "
    5  = "Middle mean draft, m: "
    6  = "Aft mean draft, m: "
    7  = "Fwd mark misplacement, m: "
    8  = "Mid mark misplacement, m: "
    9  = "Aft mark misplacement, m: "
    10 = "Apparent trim, m: "
    11 = "Fwd draft correction, m: "
    12 = "Mid draft correction, m: "
    13 = "Aft draft correction, m: "
    14 = "Fwd corrected draft, m: "
    15 = "Mid corrected draft, m: "
    16 = "Aft corrected draft, m: "
    17 = "True trim, m: "
    18 = "Deflection: "
    19 = "Mean of means corrected, m:"
    20 = "Displacement by MOMC, mt: "
    21 = "TPC, mt: "
    22 = "LCF, m: "
    23 = "First trim correction, mt:"
    24 = "MTC by MOMC: "
    25 = "MTC +: "
    26 = "MTC -: "
    27 = "MTC difference: "
    28 = "Second trim correction, mt:"
    29 = "Disp. corrected by trim, mt: "
    30 = "Constant, mt: "
    31 = "Displacement corrected, mt: "
}

$values = @{
    4  = "2.0"
    5  = "2.0"
    6  = "2.0"
    7  = "-2.095"
    8  = "1.078"
    9  = "3.067"
    10 = "0.0"
    11 = "-0.0"
    12 = "0.0"
    13 = "0.0"
    14 = "2.0"
    15 = "2.0"
    16 = "2.0"
    17 = "0.0"
    18 = "Hogging - Выгиб"
    19 = "2.0"
    20 = "2926.49"
    21 = "16.05"
    22 = "64.796"
    23 = "0.0"
    24 = "106.77"
    25 = "110.69"
    26 = "101.77"
    27 = "8.92"
    28 = "0.0"
    29 = "2926.49"
    30 = "-1142.877"
    31 = "2926.49"
}

# These values are stored as plain text (not numbers) in the source sheet,
# so format column B as Text up front to keep numeric-looking strings
# (e.g. "2.0", "-2.095") from being auto-converted to numbers.
$ws.Range("B4:B31").NumberFormat = "@"

foreach ($row in 4..31) {
    $ws.Cells.Item($row, 1).Value = $labels[$row]
    $ws.Cells.Item($row, 2).Value = $values[$row]
}
